# Continuing Issue686 [Update Features.html to match look and feel of the home
# page] - group the three existing flowchart icon clusters ("Group 13",
# "Group 29", "Group 35") together into a single group, then nudge the new
# group slightly to the right (a small manual tweak), matching the author's
# "Some more tweaks" commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Select the three top-level icon groups and group them into one shape.
$range = $s.Shapes.Range(@("Group 13", "Group 29", "Group 35"))
$newGroup = $range.Group()

# Small manual repositioning tweak applied after grouping (moves the whole
# new group to the right by ~10.08pt while leaving its vertical position and
# the relative layout of its children untouched).
$newGroup.Left = $newGroup.Left + 10.081417322834646
